$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eprepago")

# Update the data-driven test values for the new e-Prepago query
$ws.Range("D2").Value = "pruebauser01"
$ws.Range("E2").Formula = "'6789"
$ws.Range("N2").Value = "****0252"

# Make Eprepago the active/selected sheet with J10 selected
$ws.Activate() | Out-Null
$ws.Range("J10").Select() | Out-Null
